# Add 2022-Q3 data: a new "总计" summary row plus a brand-new "2022-Q3" sheet
# inserted right after "总计", pushing every later sheet one tab to the right.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" sheet: insert a new first data row (row 2) for 2022-Q3 and push
#    the existing quarters down by one row (values only, keep formatting).
#    Written as an explicit target table (A=zero-based index, B=quarter
#    label, C=holding count, D=holding value) to avoid relying on reading
#    cells back mid-script.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

$totalRows = @(
    @("2022-Q3", 14, 4.87),
    @("2022-Q2", 11, 6.01),
    @("2022-Q1", 6,  4.07),
    @("2021-Q4", 30, 6.02),
    @("2021-Q3", 15, 1.99),
    @("2021-Q2", 51, 7.6),
    @("2021-Q1", 25, 1.97),
    @("2020-Q4", 9,  0.26)
)

$rowNum = 2
$idx = 0
foreach ($row in $totalRows) {
    $wsTotal.Range("A$rowNum").Value = $idx
    $wsTotal.Range("B$rowNum").Value = $row[0]
    $wsTotal.Range("C$rowNum").Value = $row[1]
    $wsTotal.Range("D$rowNum").Value = $row[2]
    $rowNum++
    $idx++
}

# Row 9 is brand new (table used to stop at row 8) - give A9 the same
# "label" formatting (bold/bordered/centered) the rest of column A uses.
$wsTotal.Range("A8").Copy()
$wsTotal.Range("A9").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Insert a brand-new "2022-Q3" worksheet right after "总计".
# ---------------------------------------------------------------------------
$wsFirst = $wb.Worksheets.Item(1)
$wsQ3 = $wb.Worksheets.Add($null, $wsFirst)
$wsQ3.Name = "2022-Q3"

$wsQ3.Range("B1").Value = "基金代码"
$wsQ3.Range("C1").Value = "基金名称"
$wsQ3.Range("D1").Value = "基金规模"
$wsQ3.Range("E1").Value = "股票总仓位"
$wsQ3.Range("F1").Value = "仓位占比"
$wsQ3.Range("G1").Value = "持有市值(亿元)"
$wsQ3.Range("H1").Value = "仓位排名"

# B (fund code) and D/E/F/G (decimal-formatted figures) are text in the
# source data (e.g. "012930" keeps its leading zero, "47.83" stays a
# string). A bare numeric-looking literal gets auto-coerced to a number by
# this host, so prefix those with an apostrophe (Excel's quote-prefix /
# force-text convention) to keep them as text.
$q3rows = @(
    @(0,  "'012930", "中庚价值先锋股票",         "'47.83", "'94.71", "'4.56", "'2.1810", 6),
    @(1,  "'920003", "中金新锐股票A",             "'17.72", "'89.26", "'5.67", "'1.0047", 3),
    @(2,  "'000986", "太平灵活配置混合型发起式", "'12.06", "'80.26", "'5.63", "'0.6790", 8),
    @(3,  "'013004", "国泰价值领航股票A",         "'6.81",  "'92.85", "'4.17", "'0.2840", 6),
    @(4,  "'011042", "国泰价值先锋股票A",         "'4.80",  "'92.59", "'4.29", "'0.2059", 7),
    @(5,  "'920923", "中金新锐股票C",             "'3.32",  "'89.26", "'5.67", "'0.1882", 3),
    @(6,  "'001075", "宝盈转型动力灵活配置混合A", "'4.35",  "'91.90", "'3.38", "'0.1470", 6),
    @(7,  "'920002", "中金精选股票A",             "'2.95",  "'82.28", "'3.00", "'0.0885", 6),
    @(8,  "'009537", "太平行业优选股票A",         "'0.57",  "'88.65", "'7.72", "'0.0440", 4),
    @(9,  "'011043", "国泰价值先锋股票C",         "'0.37",  "'92.59", "'4.29", "'0.0159", 7),
    @(10, "'009538", "太平行业优选股票C",         "'0.17",  "'88.65", "'7.72", "'0.0131", 4),
    @(11, "'013005", "国泰价值领航股票C",         "'0.29",  "'92.85", "'4.17", "'0.0121", 6),
    @(12, "'015389", "宝盈转型动力灵活配置混合C", "'0.24",  "'91.90", "'3.38", "'0.0081", 6),
    @(13, "'920922", "中金精选股票C",             "'0.11",  "'82.28", "'3.00", "'0.0033", 6)
)

$rowNum = 2
foreach ($row in $q3rows) {
    $wsQ3.Range("A$rowNum").Value = $rowNum - 2
    $wsQ3.Range("B$rowNum").Value = $row[1]
    $wsQ3.Range("C$rowNum").Value = $row[2]
    $wsQ3.Range("D$rowNum").Value = $row[3]
    $wsQ3.Range("E$rowNum").Value = $row[4]
    $wsQ3.Range("F$rowNum").Value = $row[5]
    $wsQ3.Range("G$rowNum").Value = $row[6]
    $wsQ3.Range("H$rowNum").Value = $row[7]
    $rowNum++
}

# Match the look of every other quarter sheet: bold/bordered/centered
# "label" style on the header row and the leading index column (this is
# the same style already used by 总计!B1:D1 and 总计!A2, so reuse it instead
# of fabricating a new style entry).
$wsTotal.Range("B1:D1").Copy()
$wsQ3.Range("B1:H1").PasteSpecial(-4122)
$wsTotal.Range("A2").Copy()
$wsQ3.Range("A2:A15").PasteSpecial(-4122)

Write-Output "2022-Q3 sheet + 总计 row added"
